# Scheduled-runner style update of market/profit figures on the Ixion
# profits workbook. Each worksheet (one per Disciple of the Hand job)
# gets refreshed currentAveragePrice*/LevePrice*/LeveProfit* figures for
# a handful of leve rows, recomputed from the latest market data pull.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# --- ARM ---------------------------------------------------------------
Set-Row "ARM" 74 @{
    H = 1823.1333
    I = 1556.125
    K = 1556.125
    M = -682.125
}
Set-Row "ARM" 77 @{
    H = 1823.1333
    I = 1556.125
    K = 7780.625
    M = -3412.625
}
Set-Row "ARM" 122 @{
    H = 857473.25
    I = 918581.5600000001
    J = 1957
    K = 2755744.68
    L = 5871
    M = -2753294.68
    N = -10771
}

# --- BSM ---------------------------------------------------------------
Set-Row "BSM" 80 @{
    H = 243.51724
    J = 309.42856
    L = 309.42856
    N = -2305.42856
}
Set-Row "BSM" 83 @{
    H = 243.51724
    J = 309.42856
    L = 1547.1428
    N = -11531.1428
}
Set-Row "BSM" 122 @{
    H = 27221.111
    J = 27221.111
    L = 27221.111
    N = -37021.111
}

# --- CRP -----------------------------------------------------------------
Set-Row "CRP" 31 @{
    H = 5460.3267
    I = 1981.4839
    K = 1981.4839
    M = -1686.4839
}
Set-Row "CRP" 34 @{
    H = 5460.3267
    I = 1981.4839
    K = 1981.4839
    M = -1779.4839
}
Set-Row "CRP" 58 @{
    H = 1762.6945
    I = 1547.9
    J = 2031.1875
    K = 1547.9
    L = 2031.1875
    M = -1344.9
    N = -2437.1875
}
Set-Row "CRP" 136 @{
    H = 1762.6945
    I = 1547.9
    J = 2031.1875
    K = 4643.700000000001
    L = 6093.5625
    M = -2093.700000000001
    N = -11193.5625
}

# --- CUL -----------------------------------------------------------------
Set-Row "CUL" 12 @{
    H = 9416741
    I = 12500118
    J = 7875053
    K = 37500354
    L = 23625159
    M = -37500181
    N = -23625505
}
Set-Row "CUL" 33 @{
    H = 3575088.2
    I = 4761979
    J = 14415.714
    K = 28571874
    L = 86494.284
    M = -28571591
    N = -87060.284
}
Set-Row "CUL" 121 @{
    H = 12332.357
    I = 639.8570999999999
    J = 24024.857
    K = 1919.5713
    L = 72074.571
    M = -609.5712999999998
    N = -74694.571
}
Set-Row "CUL" 131 @{
    H = 2000993.4
    I = 11111531
    J = 1119.2927
    K = 33334593
    L = 3357.8781
    M = -33329553
    N = -13437.8781
}
Set-Row "CUL" 132 @{
    H = 2408.3635
    I = 2141.7144
    K = 19275.4296
    M = -16745.4296
}

# --- GSM -----------------------------------------------------------------
Set-Row "GSM" 80 @{
    H = 2507.1428
    I = 2512.5
    J = 2500
    K = 2512.5
    L = 2500
    M = -1514.5
    N = -4496
}
Set-Row "GSM" 83 @{
    H = 2507.1428
    I = 2512.5
    J = 2500
    K = 12562.5
    L = 12500
    M = -7570.5
    N = -22484
}
Set-Row "GSM" 122 @{
    H = 1803122.5
    I = 3243080.5
    J = 3175
    K = 9729241.5
    L = 9525
    M = -9726791.5
    N = -14425
}
Set-Row "GSM" 123 @{
    H = 21899.4
    J = 22436.875
    L = 22436.875
    N = -27336.875
}

# --- LTW -----------------------------------------------------------------
Set-Row "LTW" 68 @{
    H = 90911384
    I = 2399.875
    J = 333335330
    K = 2399.875
    L = 333335330
    M = -1650.875
    N = -333336828
}
Set-Row "LTW" 71 @{
    H = 90911384
    I = 2399.875
    J = 333335330
    K = 11999.375
    L = 1666676650
    M = -8255.375
    N = -1666684138
}
Set-Row "LTW" 122 @{
    H = 5435651.5
    I = 5502290
    J = 5002500
    K = 16506870
    L = 15007500
    M = -16504420
    N = -15012400
}

# --- WVR -----------------------------------------------------------------
Set-Row "WVR" 81 @{
    H = 1586.7273
    I = 1650.25
    J = 1417.3334
    K = 3300.5
    L = 2834.6668
    M = -2239.5
    N = -4956.6668
}
Set-Row "WVR" 84 @{
    H = 1586.7273
    I = 1650.25
    J = 1417.3334
    K = 16502.5
    L = 14173.334
    M = -11198.5
    N = -24781.334
}
Set-Row "WVR" 86 @{
    H = 21677
    J = 21677
    L = 21677
    N = -23923
}
Set-Row "WVR" 89 @{
    H = 21677
    J = 21677
    L = 108385
    N = -119617
}
Set-Row "WVR" 107 @{
    H = 52632570
    I = 90909880
    J = 1261.625
    K = 272729640
    L = 3784.875
    M = -272727720
    N = -7624.875
}
Set-Row "WVR" 122 @{
    H = 1518.4286
    I = 1031.5
    J = 2167.6667
    K = 3094.5
    L = 6503.000100000001
    M = -644.5
    N = -11403.0001
}
Set-Row "WVR" 125 @{
    H = 43715
    J = 43715
    L = 43715
    N = -53555
}
